# Refresh the cryptocurrency price ("Price", column D) and volume-change
# ("Volume(1h)", column E) figures on the active worksheet with newly
# scraped values, as produced by the scheduled GitHub Actions scraper run.
#
# Column D cells occasionally contain plain numeric-looking text (e.g.
# "354.59"). Excel's COM layer auto-converts such strings typed into
# .Value into real numbers, which would corrupt the intended text
# representation (losing trailing zeros / exact formatting). To keep
# these as literal text - matching how they were already stored in the
# workbook - we briefly force a text number format while writing the
# value, then restore the cell to the Normal style so no extra
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @(
    @{ Row = 2; Value = "51.654.32" },
    @{ Row = 3; Value = "2.917.53" },
    @{ Row = 5; Value = "354.59" },
    @{ Row = 6; Value = "110.36" },
    @{ Row = 7; Value = "0.566" },
    @{ Row = 9; Value = "0.628" },
    @{ Row = 10; Value = "39.22" },
    @{ Row = 11; Value = "0.0888" },
    @{ Row = 13; Value = "19.69" },
    @{ Row = 14; Value = "7.90" },
    @{ Row = 15; Value = "3.381.19" },
    @{ Row = 16; Value = "2.900.16" },
    @{ Row = 17; Value = "0.977" },
    @{ Row = 18; Value = "51.732.26" },
    @{ Row = 19; Value = "7.56" },
    @{ Row = 20; Value = "3.27" },
    @{ Row = 21; Value = "14.01" },
    @{ Row = 22; Value = "0.0₃0981" },
    @{ Row = 23; Value = "70.83" },
    @{ Row = 24; Value = "269.53" },
    @{ Row = 26; Value = "0.183" },
    @{ Row = 27; Value = "27.05" },
    @{ Row = 29; Value = "7.39" },
    @{ Row = 30; Value = "0.108" },
    @{ Row = 31; Value = "10.57" },
    @{ Row = 32; Value = "38.15" },
    @{ Row = 33; Value = "6.08" },
    @{ Row = 34; Value = "52.36" },
    @{ Row = 35; Value = "0.0440" },
    @{ Row = 39; Value = "18.42" },
    @{ Row = 40; Value = "2.01" },
    @{ Row = 41; Value = "2.73" },
    @{ Row = 43; Value = "22.96" },
    @{ Row = 44; Value = "122.34" },
    @{ Row = 47; Value = "3.44" },
    @{ Row = 48; Value = "2.135.62" },
    @{ Row = 49; Value = "0.252" },
    @{ Row = 51; Value = "9.09" }
)

$volumeUpdates = @(
    @{ Row = 2; Value = "  -0.95%  " },
    @{ Row = 3; Value = "  +1.14%  " },
    @{ Row = 4; Value = "  +0.05%  " },
    @{ Row = 5; Value = "  +0.63%  " },
    @{ Row = 6; Value = "  -1.43%  " },
    @{ Row = 7; Value = "  +0.61%  " },
    @{ Row = 8; Value = "  +0.01%  " },
    @{ Row = 9; Value = "  +0.99%  " },
    @{ Row = 10; Value = "  -2.48%  " },
    @{ Row = 11; Value = "  +3.31%  " },
    @{ Row = 12; Value = "  +0.79%  " },
    @{ Row = 13; Value = "  -2.12%  " },
    @{ Row = 14; Value = "  +0.79%  " },
    @{ Row = 15; Value = "  +1.20%  " },
    @{ Row = 16; Value = "  +0.02%  " },
    @{ Row = 17; Value = "  -2.12%  " },
    @{ Row = 18; Value = "  -0.82%  " },
    @{ Row = 19; Value = "  -1.87%  " },
    @{ Row = 20; Value = "  -2.74%  " },
    @{ Row = 21; Value = "  +1.81%  " },
    @{ Row = 22; Value = "  -0.22%  " },
    @{ Row = 23; Value = "  -0.42%  " },
    @{ Row = 24; Value = "  -0.34%  " },
    @{ Row = 25; Value = "  +0.99%  " },
    @{ Row = 26; Value = "  +11.38%  " },
    @{ Row = 27; Value = "  +2.51%  " },
    @{ Row = 28; Value = "  +0.10%  " },
    @{ Row = 29; Value = "  +16.11%  " },
    @{ Row = 30; Value = "  +14.90%  " },
    @{ Row = 31; Value = "  -0.04%  " },
    @{ Row = 32; Value = "  -1.72%  " },
    @{ Row = 33; Value = "  +2.85%  " },
    @{ Row = 34; Value = "  -1.72%  " },
    @{ Row = 35; Value = "  -4.36%  " },
    @{ Row = 36; Value = "  +0.05%  " },
    @{ Row = 37; Value = "  -16.10%  " },
    @{ Row = 38; Value = "  -2.86%  " },
    @{ Row = 39; Value = "  -1.65%  " },
    @{ Row = 40; Value = "  -1.01%  " },
    @{ Row = 41; Value = "  +2.81%  " },
    @{ Row = 42; Value = "  +1.19%  " },
    @{ Row = 43; Value = "  +1.49%  " },
    @{ Row = 44; Value = "  +0.61%  " },
    @{ Row = 45; Value = "  -2.61%  " },
    @{ Row = 46; Value = "  +1.33%  " },
    @{ Row = 47; Value = "  -3.73%  " },
    @{ Row = 48; Value = "  -3.14%  " },
    @{ Row = 49; Value = "  -7.25%  " },
    @{ Row = 50; Value = "  +3.24%  " },
    @{ Row = 51; Value = "  +0.17%  " }
)

foreach ($u in $priceUpdates) {
    $cell = $ws.Cells.Item($u.Row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}

foreach ($u in $volumeUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.Value
}
